# Append new scraped rows to the "ランサーズ" sheet (sheet1), matching the
# 2025-10-30 12:37:03 JST scrape run:
#   - one new job inserted at the very top (row 2)
#   - one new job inserted right before the old "non-contact sensor" row
#   - every row's "取得日時" (fetched-at) timestamp is bumped to the new run time

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-30 12:37:03"

# --- Insert the brand-new top row (AI tech-lead posting) -------------------
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆効率化"
$ws.Range("F2").Style = "Hyperlink"

# --- Insert the brand-new row just before the sensor-idea posting ----------
# (after the first insert, the old rows 2-9 now live at rows 3-10, so the
#  old row 8 "non-contact sensor" posting currently sits at row 9)
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = $newTimestamp
$ws.Range("B9").Value = "評価基板設計・製造について"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5423728"
$ws.Range("G9").Value = 18

# --- Refresh the timestamp on every other (pre-existing) row ---------------
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp

# --- Hyperlinks for the two rows that now extend past the old table end ----
# (rows 10 and 11 are the old "sensor idea" and "Google Workspace" postings,
#  shifted down past the sheet's former last row; give them live hyperlinks
#  just like every other URL cell in column F)
$ws.Hyperlinks.Add($ws.Range("F10"), $ws.Range("F10").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), $ws.Range("F11").Value2) | Out-Null
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("F11").Style = "Hyperlink"
